# Updates the cryptocurrency price/volume table (and reorders the
# Aptos / Fetch.AI rows) to match the latest scrape, per the commit
# "Updated cryptos list ... with GitHub Actions".
#
# Each entry below is one cell whose text content changed. Cells whose
# new text parses as a plain number (e.g. "0.997", "522.02") are written
# with a leading apostrophe so Excel keeps storing them as text (matching
# the original "Price"/"Volume" columns, which are text, not numeric,
# cells) and the cell's style is then reset to "Normal" so no stray
# number-format/quote-prefix formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '58.532.87'; Numeric = $false },
    @{ Cell = 'E2'; Value = '  -0.85%  '; Numeric = $false },
    @{ Cell = 'D3'; Value = '2.474.25'; Numeric = $false },
    @{ Cell = 'E3'; Value = '  -1.15%  '; Numeric = $false },
    @{ Cell = 'D4'; Value = '0.997'; Numeric = $true },
    @{ Cell = 'E4'; Value = '  -0.27%  '; Numeric = $false },
    @{ Cell = 'D5'; Value = '522.02'; Numeric = $true },
    @{ Cell = 'E5'; Value = '  -2.32%  '; Numeric = $false },
    @{ Cell = 'D6'; Value = '134.56'; Numeric = $true },
    @{ Cell = 'E6'; Value = '  -1.32%  '; Numeric = $false },
    @{ Cell = 'D7'; Value = '0.997'; Numeric = $true },
    @{ Cell = 'E7'; Value = '  -0.20%  '; Numeric = $false },
    @{ Cell = 'D8'; Value = '0.559'; Numeric = $true },
    @{ Cell = 'E8'; Value = '  -1.32%  '; Numeric = $false },
    @{ Cell = 'D9'; Value = '2.478.92'; Numeric = $false },
    @{ Cell = 'E9'; Value = '  -1.04%  '; Numeric = $false },
    @{ Cell = 'D10'; Value = '0.0984'; Numeric = $true },
    @{ Cell = 'E10'; Value = '  -3.53%  '; Numeric = $false },
    @{ Cell = 'E11'; Value = '  -0.86%  '; Numeric = $false },
    @{ Cell = 'D12'; Value = '5.34'; Numeric = $true },
    @{ Cell = 'E12'; Value = '  -0.90%  '; Numeric = $false },
    @{ Cell = 'D13'; Value = '0.339'; Numeric = $true },
    @{ Cell = 'E13'; Value = '  -2.71%  '; Numeric = $false },
    @{ Cell = 'D14'; Value = '2.907.03'; Numeric = $false },
    @{ Cell = 'E14'; Value = '  -1.32%  '; Numeric = $false },
    @{ Cell = 'D15'; Value = '58.350.24'; Numeric = $false },
    @{ Cell = 'E15'; Value = '  -0.99%  '; Numeric = $false },
    @{ Cell = 'D16'; Value = '22.23'; Numeric = $true },
    @{ Cell = 'E16'; Value = '  -2.21%  '; Numeric = $false },
    @{ Cell = 'D17'; Value = '0.0000135'; Numeric = $true },
    @{ Cell = 'E17'; Value = '  -2.27%  '; Numeric = $false },
    @{ Cell = 'D18'; Value = '2.476.99'; Numeric = $false },
    @{ Cell = 'E18'; Value = '  -1.22%  '; Numeric = $false },
    @{ Cell = 'D19'; Value = '10.69'; Numeric = $true },
    @{ Cell = 'E19'; Value = '  -3.18%  '; Numeric = $false },
    @{ Cell = 'D20'; Value = '321.49'; Numeric = $true },
    @{ Cell = 'E20'; Value = '  -0.65%  '; Numeric = $false },
    @{ Cell = 'D21'; Value = '4.19'; Numeric = $true },
    @{ Cell = 'E21'; Value = '  -1.84%  '; Numeric = $false },
    @{ Cell = 'E22'; Value = '  +0.01%  '; Numeric = $false },
    @{ Cell = 'D23'; Value = '5.74'; Numeric = $true },
    @{ Cell = 'E23'; Value = '  -5.52%  '; Numeric = $false },
    @{ Cell = 'D24'; Value = '64.78'; Numeric = $true },
    @{ Cell = 'E24'; Value = '  -0.68%  '; Numeric = $false },
    @{ Cell = 'D25'; Value = '0.410'; Numeric = $true },
    @{ Cell = 'E25'; Value = '  -2.83%  '; Numeric = $false },
    @{ Cell = 'D26'; Value = '0.996'; Numeric = $true },
    @{ Cell = 'E26'; Value = '  -0.29%  '; Numeric = $false },
    @{ Cell = 'D27'; Value = '0.161'; Numeric = $true },
    @{ Cell = 'E27'; Value = '  -2.03%  '; Numeric = $false },
    @{ Cell = 'D28'; Value = '7.37'; Numeric = $true },
    @{ Cell = 'E28'; Value = '  -2.36%  '; Numeric = $false },
    @{ Cell = 'D29'; Value = '0.0₃0751'; Numeric = $false },
    @{ Cell = 'E29'; Value = '  -1.70%  '; Numeric = $false },
    @{ Cell = 'D30'; Value = '169.12'; Numeric = $true },
    @{ Cell = 'E30'; Value = '  -0.63%  '; Numeric = $false },
    @{ Cell = 'D31'; Value = '1.70'; Numeric = $true },
    @{ Cell = 'E31'; Value = '  -2.86%  '; Numeric = $false },
    @{ Cell = 'B32'; Value = 'Fetch.AI'; Numeric = $false },
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; Numeric = $false },
    @{ Cell = 'D32'; Value = '1.18'; Numeric = $true },
    @{ Cell = 'E32'; Value = '  +0.94%  '; Numeric = $false },
    @{ Cell = 'B33'; Value = 'Aptos'; Numeric = $false },
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Numeric = $false },
    @{ Cell = 'D33'; Value = '6.29'; Numeric = $true },
    @{ Cell = 'E33'; Value = '  -3.35%  '; Numeric = $false },
    @{ Cell = 'E34'; Value = '  -0.05%  '; Numeric = $false },
    @{ Cell = 'D35'; Value = '0.998'; Numeric = $true },
    @{ Cell = 'E35'; Value = '  -0.07%  '; Numeric = $false },
    @{ Cell = 'D36'; Value = '18.12'; Numeric = $true },
    @{ Cell = 'E36'; Value = '  -1.39%  '; Numeric = $false },
    @{ Cell = 'D37'; Value = '1.34'; Numeric = $true },
    @{ Cell = 'E37'; Value = '  -1.50%  '; Numeric = $false },
    @{ Cell = 'D38'; Value = '4.00'; Numeric = $true },
    @{ Cell = 'E38'; Value = '  -1.29%  '; Numeric = $false },
    @{ Cell = 'D39'; Value = '36.53'; Numeric = $true },
    @{ Cell = 'E39'; Value = '  -0.54%  '; Numeric = $false },
    @{ Cell = 'D40'; Value = '1.48'; Numeric = $true },
    @{ Cell = 'E40'; Value = '  -3.26%  '; Numeric = $false },
    @{ Cell = 'D41'; Value = '0.802'; Numeric = $true },
    @{ Cell = 'E41'; Value = '  -0.09%  '; Numeric = $false },
    @{ Cell = 'D42'; Value = '3.46'; Numeric = $true },
    @{ Cell = 'E42'; Value = '  -3.43%  '; Numeric = $false },
    @{ Cell = 'D43'; Value = '274.88'; Numeric = $true },
    @{ Cell = 'E43'; Value = '  -2.80%  '; Numeric = $false },
    @{ Cell = 'D44'; Value = '5.08'; Numeric = $true },
    @{ Cell = 'E44'; Value = '  +0.85%  '; Numeric = $false },
    @{ Cell = 'D45'; Value = '0.598'; Numeric = $true },
    @{ Cell = 'E45'; Value = '  -0.30%  '; Numeric = $false },
    @{ Cell = 'D46'; Value = '124.72'; Numeric = $true },
    @{ Cell = 'E46'; Value = '  -4.44%  '; Numeric = $false },
    @{ Cell = 'D47'; Value = '0.0911'; Numeric = $true },
    @{ Cell = 'E47'; Value = '  -1.49%  '; Numeric = $false },
    @{ Cell = 'D48'; Value = '0.0490'; Numeric = $true },
    @{ Cell = 'E48'; Value = '  -2.14%  '; Numeric = $false },
    @{ Cell = 'D49'; Value = '0.0213'; Numeric = $true },
    @{ Cell = 'E49'; Value = '  -2.28%  '; Numeric = $false },
    @{ Cell = 'D50'; Value = '17.08'; Numeric = $true },
    @{ Cell = 'E50'; Value = '  -1.43%  '; Numeric = $false },
    @{ Cell = 'D51'; Value = '1.737.46'; Numeric = $false },
    @{ Cell = 'E51'; Value = '  -1.12%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $cell.Value = "'" + $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
